# #5: cash & deposit done
# Rebuild the "存款" (deposit) sheet: fix the header row (which had stray
# data instead of labels) and add the bank / deposit_type / currency
# columns plus the standard metadata columns (property_category .. index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Header row (row 1) -------------------------------------------------
# Was: B1=<bank name>, C1=<deposit type>, D1=<currency>, E1=<legislator>, F1=<amount>
# Now: proper column headers, extended through M1.
$ws.Range("G1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

$ws.Range("B1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Data rows (2-7) -----------------------------------------------------
# Columns A (index), E (owner) and F (amount) already hold the right data;
# B/C/D get corrected bank / deposit_type / currency text, and G:M are new.

$banks = @("玉山商業銀行北新分行", "玉山商業銀行北新分行", "中國信託商業銀行板橋分行", "華南商業銀行華江分行", "華南商業銀行文山分行", "華南商業纟1行文山另:行")

for ($i = 0; $i -lt 6; $i++) {
    $r = 2 + $i

    $ws.Range("B$r").Value = $banks[$i]
    $ws.Range("C$r").Value = "活期儲蓄存款"
    $ws.Range("D$r").Value = "新臺幣"

    $ws.Range("G1").Copy()
    $ws.Range("G$r" + ":M$r").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"

    # Force the date to be stored as literal text, not auto-converted to a
    # date serial, then strip the number-format override again so the cell
    # ends up with the same (default) style as its neighbours.
    $ws.Range("I$r").NumberFormat = "@"
    $ws.Range("I$r").Value = "2011-12-31"
    $ws.Range("H$r").Copy()
    $ws.Range("I$r").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("J$r").Value = "羅明才"
    $ws.Range("K$r").Value = 879
    $ws.Range("L$r").Value = "tmp94331"
    $ws.Range("M$r").Value = $ws.Range("A$r").Value2
}
